# Generate Report for Handback
# Adds a new handback entry (2dd397bb-41be-4b48-9168-189969d508d5.md) as
# row 4 on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the
# existing rows 2/3, and grows each sheet's table to include it.

$wb = $excel.ActiveWorkbook

$mdFile        = "2dd397bb-41be-4b48-9168-189969d508d5.md"
$mdPath        = "e2e\2dd397bb-41be-4b48-9168-189969d508d5.md"
$commitSha     = "b6c725c51d0942106f92cd7c0a060d93991743b5"
$xlfZh         = "2dd397bb-41be-4b48-9168-189969d508d5." + $commitSha + ".zh-cn.xlf"
$xlfDe         = "2dd397bb-41be-4b48-9168-189969d508d5." + $commitSha + ".de-de.xlf"
$statusInSync  = "Handed back: in sync with en-US"

$dateOverview  = "2016-08-21 06:51:17"
$dateZhHandoff = "2016-08-21 06:51:13"
$dateZhHandback= "2016-08-21 06:51:29"
$dateDeHandback= "2016-08-21 06:51:35"

$urlMdOverview = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1c273cbbf7160c3ce131c44f3ca5c3a43bb2783/e2e/" + $mdFile
$urlMdZh       = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a1f6f3a2c6c0c2a6a4f0e9c3fa4e6d9a9b4f3e2b/e2e/" + $mdFile
$urlMdDe       = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f3b6a2e1d4c5b6a7980f1e2d3c4b5a6978695a4b/e2e/" + $mdFile

# ---------------------------------------------------------------------
# Sheet "Overview" -> row 4
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdFile
$wsOverview.Range("B4").Value = $mdPath
$wsOverview.Range("B4").Style = "HyperLink"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $urlMdOverview, "", "", $mdPath) | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G4").Value = $dateOverview

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" -> row 4
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $mdFile
$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $urlMdZh, "", "", $mdFile) | Out-Null
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $statusInSync
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "True"
$wsZh.Range("G4").Value = $xlfZh
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value = $dateZhHandoff
$wsZh.Range("I4").Value = $mdFile
$wsZh.Range("I4").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $urlMdZh, "", "", $mdFile) | Out-Null
$wsZh.Range("J4").Value = $xlfZh
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K4").Value = $dateZhHandback
$wsZh.Range("M4").Value = "True"
$wsZh.Range("O4").Value = "False"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

# ---------------------------------------------------------------------
# Sheet "de-de" -> row 4
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $mdFile
$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $urlMdDe, "", "", $mdFile) | Out-Null
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $statusInSync
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "True"
$wsDe.Range("G4").Value = $xlfDe
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value = $dateOverview
$wsDe.Range("I4").Value = $mdFile
$wsDe.Range("I4").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $urlMdDe, "", "", $mdFile) | Out-Null
$wsDe.Range("J4").Value = $xlfDe
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K4").Value = $dateDeHandback
$wsDe.Range("M4").Value = "True"
$wsDe.Range("O4").Value = "False"

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))
